$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "factorielTest/testComplexNumbers"
$ws.Range("B10").Value = $true
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = 0.0073391999999999997

$ws.Range("E2").Value = 0.010650199999999999
$ws.Range("E3").Value = 0.0029557999999999997
$ws.Range("E4").Value = 0.0028760999999999999
$ws.Range("E5").Value = 0.0067183000000000008
$ws.Range("E6").Value = 0.0067672999999999995
$ws.Range("E7").Value = 0.0028725999999999999
$ws.Range("E8").Value = 0.0025345000000000003
$ws.Range("E9").Value = 0.0032951000000000005
